$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("B2").Value = 1518

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 1377

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 521

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 426

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 413

# Copy the style from A3 (existing styled cell) to the newly added A4:A6 cells
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0
